# #5: insurance, claim, debt, investment done
# Extends the "債務" (debt) sheet with property_category/category/date/
# legislator_name/legislator_id/source_file/index columns (H:N), matching
# the pattern already used on the other sheets, and rewrites the header
# row (B1:G1) to use the real field names instead of stray sample data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("債務")

$xlPasteFormats = -4122

# --- extend formatting for the new columns (H:N) ------------------------
# Header row (row 1) uses the bordered/bold style already on G1.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:N1").PasteSpecial($xlPasteFormats) | Out-Null

# Data rows (2-6) use the plain style already on column G of that row.
for ($r = 2; $r -le 6; $r++) {
    $ws.Range("G$r").Copy() | Out-Null
    $ws.Range("H" + $r + ":N" + $r).PasteSpecial($xlPasteFormats) | Out-Null
}

# --- header row: real field names ---------------------------------------
$ws.Cells.Item(1, 2).Value = "species"
$ws.Cells.Item(1, 3).Value = "debtor"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "total"
$ws.Cells.Item(1, 6).Value = "register_date"
$ws.Cells.Item(1, 7).Value = "register_reason"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- data rows: fill in the new metadata columns H:N ---------------------
$indices = @{ 2 = 108; 3 = 110; 4 = 111; 5 = 112; 6 = 113 }

foreach ($r in 2..6) {
    $ws.Cells.Item($r, 8).Value = "debt"
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 10).Value = "2012-05-01"
    $ws.Cells.Item($r, 11).Value = "管碧玲"
    $ws.Cells.Item($r, 12).Value = 1374
    $ws.Cells.Item($r, 13).Value = "tmpf0df1"
    $ws.Cells.Item($r, 14).Value = $indices[$r]
}

Write-Host "sheet5 (債務) updated"
